# Generate Report for Archive
# - Flip the "Ready for handoff" status to "In Translation" everywhere it
#   appears (Overview!E2:F2, zh-cn!C2, de-de!C2 all share that string).
# - Narrow the now-shorter "Status" columns (Overview E:F, zh-cn C, de-de C)
#   to match the refreshed report layout.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($row = 1; $row -le $used.Rows.Count; $row++) {
        for ($col = 1; $col -le $used.Columns.Count; $col++) {
            $cell = $used.Cells.Item($row, $col)
            $val = $cell.Value()
            # NB: keep the string literal on the left of -eq -- PowerShell's
            # comparison operator coerces the right-hand side to the type of
            # the left-hand side, so "$boolCell.Value -eq $oldStatus" would
            # otherwise happily (and wrongly) match boolean TRUE cells too.
            if ($oldStatus -eq $val) {
                $cell.Value = $newStatus
                $cell.EntireColumn.ColumnWidth = 12.5
            }
        }
    }
}
